# Weekly update: insert a new price record as the new first data row (row 50)
# for "Agrícola del Norte S.A. de Arica" - Ají, pushing the existing rows
# 50-84 down by one (they become rows 51-85).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 50; this shifts rows 50:84 down to 51:85
# and extends the used range to A1:R85 automatically.
$ws.Rows.Item(50).Insert()

# Populate the newly inserted row 50 with the new weekly record.
$ws.Cells.Item(50, 1).Value = 1
$ws.Cells.Item(50, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(50, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(50, 4).Value = 44729
$ws.Cells.Item(50, 5).Value = 15
$ws.Cells.Item(50, 6).Value = 100112021
$ws.Cells.Item(50, 7).Value = "Ají"
$ws.Cells.Item(50, 8).Value = "Inferno"
$ws.Cells.Item(50, 9).Value = "Primera"
$ws.Cells.Item(50, 10).Value = 130
$ws.Cells.Item(50, 11).Value = 16000
$ws.Cells.Item(50, 12).Value = 17000
$ws.Cells.Item(50, 13).Value = 16500
$ws.Cells.Item(50, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(50, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(50, 16).Value = 1100
$ws.Cells.Item(50, 17).Value = 15
$ws.Cells.Item(50, 18).Value = "Hortaliza"
